$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 11
$ws.Cells.Item($row, 1).Value = "'2025-01-06"
$ws.Cells.Item($row, 2).Value = "22:28:14"
$ws.Cells.Item($row, 3).Value = "Monday"
$ws.Cells.Item($row, 4).Value = "'01"
$ws.Cells.Item($row, 5).Value = 127448
$ws.Cells.Item($row, 6).Value = 143614
$ws.Cells.Item($row, 7).Value = 169023
$ws.Cells.Item($row, 8).Value = 158495
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142137
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192470
$ws.Cells.Item($row, 14).Value = 114981
$ws.Cells.Item($row, 15).Value = 45612
$ws.Cells.Item($row, 16).Value = 28331
$ws.Cells.Item($row, 17).Value = 64103
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 47957
$ws.Cells.Item($row, 20).Value = -1
